$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.559.62'
$ws.Range("E2").Value = '  +1.50%  '

$ws.Range("D3").Value = '3.378.26'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.77'
$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = '  +5.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '188.10'
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = '  -1.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.595'
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E8").Value = '  +1.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.183'
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.588'
$ws.Range("C10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E10").Value = '  +0.26%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.59'
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E11").Value = '  +0.97%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000275'
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("D13").Value = '3.926.27'
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '641.33'
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = '  +6.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.62'
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = '  -1.30%  '

$ws.Range("D16").Value = '67.695.52'
$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("D17").Value = '3.388.33'
$ws.Range("E17").Value = '  +0.64%  '

$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.05'
$ws.Range("C19").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.19'
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.911'
$ws.Range("C21").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = '  +0.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.07'
$ws.Range("C22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = '  -1.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.12'
$ws.Range("C23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = '  +1.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '99.38'
$ws.Range("C24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = '  -1.10%  '

$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.71'
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = '  +1.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.41'
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = '  +5.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.70'
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.96'
$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = '  +3.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '610.18'
$ws.Range("C31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = '  +3.99%  '

$ws.Range("E32").Value = '  -2.59%  '

$ws.Range("D33").Value = '3.998.39'
$ws.Range("E33").Value = '  +7.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.13'
$ws.Range("C34").Copy()
$ws.Range("D34").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("E35").Value = '  +0.90%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.18'
$ws.Range("C37").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.85'
$ws.Range("C38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = '  +6.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.132'
$ws.Range("C39").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E39").Value = '  +3.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.83'
$ws.Range("C40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = '  -0.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.26'
$ws.Range("C41").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").Value = '0.0₃0707'
$ws.Range("E42").Value = '  -0.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.345'
$ws.Range("C43").Copy()
$ws.Range("D43").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("C44").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("E45").Value = '  +0.50%  '

$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.59'
$ws.Range("C47").Copy()
$ws.Range("D47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E47").Value = '  -0.18%  '

$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("E49").Value = '  +10.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.82'
$ws.Range("C50").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E50").Value = '  -21.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.50'
$ws.Range("C51").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E51").Value = '  +3.59%  '
